$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 270.16666
$ws.Range("I18").Value = 270.16666
$ws.Range("K18").Value = 270.16666
$ws.Range("M18").Value = 13.83334000000002
$ws.Range("H28").Value = 2570
$ws.Range("I28").Value = 1663.1666
$ws.Range("K28").Value = 1663.1666
$ws.Range("M28").Value = -1178.1666
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("H76").Value = 3754.5557
$ws.Range("I76").Value = 3698.8572
$ws.Range("J76").Value = 3949.5
$ws.Range("K76").Value = 3698.8572
$ws.Range("L76").Value = 3949.5
$ws.Range("M76").Value = -3383.8572
$ws.Range("N76").Value = -4579.5
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H79").Value = 3754.5557
$ws.Range("I79").Value = 3698.8572
$ws.Range("J79").Value = 3949.5
$ws.Range("K79").Value = 3698.8572
$ws.Range("L79").Value = 3949.5
$ws.Range("M79").Value = -2606.8572
$ws.Range("N79").Value = -6133.5
$ws.Range("H87").Value = 83840.92
$ws.Range("J87").Value = 83840.92
$ws.Range("L87").Value = 83840.92
$ws.Range("N87").Value = -86336.92
$ws.Range("H90").Value = 83840.92
$ws.Range("J90").Value = 83840.92
$ws.Range("L90").Value = 251522.76
$ws.Range("N90").Value = -264002.76
$ws.Range("H98").Value = 1584.7059
$ws.Range("I98").Value = 1584.7059
$ws.Range("K98").Value = 1584.7059
$ws.Range("M98").Value = -86.70589999999993
$ws.Range("H100").Value = 45365.78
$ws.Range("I100").Value = 48829.19
$ws.Range("K100").Value = 48829.19
$ws.Range("M100").Value = -48288.19
$ws.Range("H106").Value = 18491.166
$ws.Range("I106").Value = 9227.959999999999
$ws.Range("K106").Value = 9227.959999999999
$ws.Range("M106").Value = -8596.959999999999
$ws.Range("H107").Value = 22491.5
$ws.Range("I107").Value = 9487.723
$ws.Range("K107").Value = 9487.723
$ws.Range("M107").Value = -7567.723
$ws.Range("H110").Value = 44641.668
$ws.Range("J110").Value = 44641.668
$ws.Range("L110").Value = 44641.668
$ws.Range("N110").Value = -52821.668
$ws.Range("H111").Value = 3452.2144
$ws.Range("I111").Value = 3194.3333
$ws.Range("J111").Value = 4999.5
$ws.Range("K111").Value = 9582.999899999999
$ws.Range("L111").Value = 14998.5
$ws.Range("M111").Value = -6515.999899999999
$ws.Range("N111").Value = -21132.5
$ws.Range("H116").Value = 23496.834
$ws.Range("I116").Value = 9999
$ws.Range("K116").Value = 9999
$ws.Range("H121").Value = 2884.4
$ws.Range("J121").Value = 3289.5833
$ws.Range("L121").Value = 9868.749899999999
$ws.Range("N121").Value = -13362.7499
$ws.Range("H122").Value = 1584.7059
$ws.Range("I122").Value = 1584.7059
$ws.Range("K122").Value = 4754.1177
$ws.Range("M122").Value = -2304.1177
$ws.Range("H125").Value = 748.5
$ws.Range("J125").Value = 748.5
$ws.Range("L125").Value = 6736.5
$ws.Range("N125").Value = -11656.5
$ws.Range("H135").Value = 4568.3887
$ws.Range("I135").Value = 5225.769
$ws.Range("K135").Value = 47031.921
$ws.Range("M135").Value = -44496.921
$ws.Range("H137").Value = 7447.8423
$ws.Range("I137").Value = 2937.353
$ws.Range("J137").Value = 14115.521
$ws.Range("K137").Value = 8812.059000000001
$ws.Range("L137").Value = 42346.563
$ws.Range("M137").Value = -6262.059000000001
$ws.Range("N137").Value = -47446.563
$ws.Range("H138").Value = 6249.7036
$ws.Range("J138").Value = 7719.875
$ws.Range("L138").Value = 23159.625
$ws.Range("N138").Value = -33439.625
$ws.Range("H141").Value = 2454.6924
$ws.Range("I141").Value = 2422.5454
$ws.Range("K141").Value = 7267.6362
$ws.Range("M141").Value = -2087.6362
$ws.Range("M116").Value = -6557
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("N78").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3808.898
$ws.Range("I2").Value = 3724.7144
$ws.Range("J2").Value = 4019.3572
$ws.Range("K2").Value = 3724.7144
$ws.Range("L2").Value = 4019.3572
$ws.Range("M2").Value = -3611.7144
$ws.Range("N2").Value = -4245.3572
$ws.Range("H32").Value = 3938.8523
$ws.Range("I32").Value = 3704.8953
$ws.Range("J32").Value = 13999
$ws.Range("K32").Value = 3704.8953
$ws.Range("L32").Value = 13999
$ws.Range("M32").Value = -3417.8953
$ws.Range("H35").Value = 22653.666
$ws.Range("J35").Value = 29999
$ws.Range("L35").Value = 29999
$ws.Range("H61").Value = 6035.6387
$ws.Range("I61").Value = 4220.1724
$ws.Range("K61").Value = 4220.1724
$ws.Range("M61").Value = -4008.1724
$ws.Range("H74").Value = 8980.666999999999
$ws.Range("I74").Value = 998.2
$ws.Range("J74").Value = 14682.429
$ws.Range("K74").Value = 998.2
$ws.Range("L74").Value = 14682.429
$ws.Range("M74").Value = -124.2
$ws.Range("N74").Value = -16430.429
$ws.Range("H77").Value = 8980.666999999999
$ws.Range("I77").Value = 998.2
$ws.Range("J77").Value = 14682.429
$ws.Range("K77").Value = 4991
$ws.Range("L77").Value = 73412.145
$ws.Range("M77").Value = -623
$ws.Range("N77").Value = -82148.145
$ws.Range("H96").Value = 97669
$ws.Range("J96").Value = 97669
$ws.Range("L96").Value = 97669
$ws.Range("N96").Value = -103161
$ws.Range("H101").Value = 57371.57
$ws.Range("J101").Value = 57371.57
$ws.Range("L101").Value = 57371.57
$ws.Range("N101").Value = -63861.57
$ws.Range("H105").Value = 29799
$ws.Range("J105").Value = 29799
$ws.Range("L105").Value = 29799
$ws.Range("H110").Value = 2287.5557
$ws.Range("I110").Value = 2373.5
$ws.Range("J110").Value = 1600
$ws.Range("K110").Value = 2373.5
$ws.Range("L110").Value = 1600
$ws.Range("M110").Value = -328.5
$ws.Range("N110").Value = -5690
$ws.Range("H113").Value = 45000
$ws.Range("J113").Value = 45000
$ws.Range("L113").Value = 45000
$ws.Range("H116").Value = 3808.898
$ws.Range("I116").Value = 3724.7144
$ws.Range("J116").Value = 4019.3572
$ws.Range("K116").Value = 3724.7144
$ws.Range("L116").Value = 4019.3572
$ws.Range("M116").Value = -1430.7144
$ws.Range("N116").Value = -8607.3572
$ws.Range("H122").Value = 2100.7693
$ws.Range("I122").Value = 1940.5
$ws.Range("K122").Value = 5821.5
$ws.Range("M122").Value = -3371.5
$ws.Range("H132").Value = 1975.8572
$ws.Range("I132").Value = 1923.6842
$ws.Range("K132").Value = 5771.0526
$ws.Range("M132").Value = -3241.0526
$ws.Range("H136").Value = 6035.6387
$ws.Range("I136").Value = 4220.1724
$ws.Range("K136").Value = 12660.5172
$ws.Range("M136").Value = -10110.5172
$ws.Range("N32").Value = -14573
$ws.Range("N35").Value = -30811
$ws.Range("N105").Value = -36787
$ws.Range("N113").Value = -53678

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3808.898
$ws.Range("I3").Value = 3724.7144
$ws.Range("J3").Value = 4019.3572
$ws.Range("K3").Value = 3724.7144
$ws.Range("L3").Value = 4019.3572
$ws.Range("M3").Value = -3610.7144
$ws.Range("N3").Value = -4247.3572
$ws.Range("H20").Value = 22071.096
$ws.Range("I20").Value = 24433.2
$ws.Range("J20").Value = 16165.833
$ws.Range("K20").Value = 24433.2
$ws.Range("L20").Value = 16165.833
$ws.Range("M20").Value = -24186.2
$ws.Range("N20").Value = -16659.833
$ws.Range("H37").Value = 662.8
$ws.Range("I37").Value = 436.55554
$ws.Range("K37").Value = 436.55554
$ws.Range("M37").Value = -299.55554
$ws.Range("H64").Value = 335
$ws.Range("I64").Value = 481.66666
$ws.Range("J64").Value = 286.1111
$ws.Range("K64").Value = 481.66666
$ws.Range("L64").Value = 286.1111
$ws.Range("M64").Value = -256.66666
$ws.Range("N64").Value = -736.1111000000001
$ws.Range("H67").Value = 335
$ws.Range("I67").Value = 481.66666
$ws.Range("J67").Value = 286.1111
$ws.Range("K67").Value = 481.66666
$ws.Range("L67").Value = 286.1111
$ws.Range("M67").Value = 298.33334
$ws.Range("N67").Value = -1846.1111
$ws.Range("H80").Value = 313.9
$ws.Range("I80").Value = 241.1
$ws.Range("J80").Value = 386.7
$ws.Range("K80").Value = 241.1
$ws.Range("L80").Value = 386.7
$ws.Range("M80").Value = 756.9
$ws.Range("N80").Value = -2382.7
$ws.Range("H82").Value = 51525.438
$ws.Range("J82").Value = 91244.75
$ws.Range("L82").Value = 91244.75
$ws.Range("N82").Value = -92010.75
$ws.Range("H83").Value = 313.9
$ws.Range("I83").Value = 241.1
$ws.Range("J83").Value = 386.7
$ws.Range("K83").Value = 1205.5
$ws.Range("L83").Value = 1933.5
$ws.Range("M83").Value = 3786.5
$ws.Range("N83").Value = -11917.5
$ws.Range("H85").Value = 51525.438
$ws.Range("J85").Value = 91244.75
$ws.Range("L85").Value = 91244.75
$ws.Range("N85").Value = -93896.75
$ws.Range("H86").Value = 348170.1
$ws.Range("I86").Value = 716837.7
$ws.Range("K86").Value = 716837.7
$ws.Range("M86").Value = -715714.7
$ws.Range("H89").Value = 348170.1
$ws.Range("I89").Value = 716837.7
$ws.Range("K89").Value = 3584188.5
$ws.Range("M89").Value = -3578572.5
$ws.Range("H102").Value = 14863.875
$ws.Range("I102").Value = 14863.875
$ws.Range("K102").Value = 14863.875
$ws.Range("M102").Value = -11618.875
$ws.Range("H105").Value = 4967
$ws.Range("I105").Value = 4668.4287
$ws.Range("J105").Value = 5803
$ws.Range("K105").Value = 4668.4287
$ws.Range("L105").Value = 5803
$ws.Range("M105").Value = -2921.4287
$ws.Range("N105").Value = -9297
$ws.Range("H107").Value = 4603.4
$ws.Range("I107").Value = 5725
$ws.Range("J107").Value = 3232.5557
$ws.Range("K107").Value = 5725
$ws.Range("L107").Value = 3232.5557
$ws.Range("M107").Value = -3805
$ws.Range("N107").Value = -7072.5557
$ws.Range("H130").Value = 68569.14
$ws.Range("J130").Value = 68569.14
$ws.Range("L130").Value = 68569.14
$ws.Range("N130").Value = -78609.14
$ws.Range("H132").Value = 84943.3
$ws.Range("J132").Value = 84943.3
$ws.Range("L132").Value = 84943.3
$ws.Range("N132").Value = -95063.3
$ws.Range("H134").Value = 10610.462
$ws.Range("I134").Value = 5482.4053
$ws.Range("K134").Value = 16447.2159
$ws.Range("M134").Value = -13912.2159

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 82049560
$ws.Range("I4").Value = 50000
$ws.Range("K4").Value = 50000
$ws.Range("H31").Value = 3273.2354
$ws.Range("I31").Value = 2044.6666
$ws.Range("J31").Value = 4655.375
$ws.Range("K31").Value = 2044.6666
$ws.Range("L31").Value = 4655.375
$ws.Range("M31").Value = -1749.6666
$ws.Range("N31").Value = -5245.375
$ws.Range("H34").Value = 3273.2354
$ws.Range("I34").Value = 2044.6666
$ws.Range("J34").Value = 4655.375
$ws.Range("K34").Value = 2044.6666
$ws.Range("L34").Value = 4655.375
$ws.Range("M34").Value = -1842.6666
$ws.Range("N34").Value = -5059.375
$ws.Range("H41").Value = 27761.883
$ws.Range("J41").Value = 42216.89
$ws.Range("L41").Value = 42216.89
$ws.Range("N41").Value = -43072.89
$ws.Range("H50").Value = 40244.625
$ws.Range("J50").Value = 40244.625
$ws.Range("L50").Value = 40244.625
$ws.Range("N50").Value = -41494.625
$ws.Range("H51").Value = 43328
$ws.Range("J51").Value = 43328
$ws.Range("L51").Value = 43328
$ws.Range("N51").Value = -44800
$ws.Range("H59").Value = 70894.89999999999
$ws.Range("J59").Value = 69883.44500000001
$ws.Range("L59").Value = 69883.44500000001
$ws.Range("N59").Value = -72173.44500000001
$ws.Range("H60").Value = 35621.125
$ws.Range("I60").Value = 15000
$ws.Range("J60").Value = 38567
$ws.Range("K60").Value = 15000
$ws.Range("L60").Value = 38567
$ws.Range("M60").Value = -14489
$ws.Range("N60").Value = -39589
$ws.Range("H61").Value = 43328
$ws.Range("J61").Value = 43328
$ws.Range("L61").Value = 43328
$ws.Range("N61").Value = -44024
$ws.Range("H86").Value = 2879.0588
$ws.Range("I86").Value = 2449.1
$ws.Range("J86").Value = 3493.2856
$ws.Range("K86").Value = 2449.1
$ws.Range("L86").Value = 3493.2856
$ws.Range("M86").Value = -1326.1
$ws.Range("N86").Value = -5739.2856
$ws.Range("H87").Value = 64990.332
$ws.Range("J87").Value = 64990.332
$ws.Range("L87").Value = 64990.332
$ws.Range("N87").Value = -67362.33199999999
$ws.Range("H88").Value = 7752.5713
$ws.Range("J88").Value = 7752.5713
$ws.Range("L88").Value = 7752.5713
$ws.Range("N88").Value = -8564.5713
$ws.Range("H89").Value = 2879.0588
$ws.Range("I89").Value = 2449.1
$ws.Range("J89").Value = 3493.2856
$ws.Range("K89").Value = 12245.5
$ws.Range("L89").Value = 17466.428
$ws.Range("M89").Value = -6629.5
$ws.Range("N89").Value = -28698.428
$ws.Range("H90").Value = 64990.332
$ws.Range("J90").Value = 64990.332
$ws.Range("L90").Value = 194970.996
$ws.Range("N90").Value = -206826.996
$ws.Range("H91").Value = 7752.5713
$ws.Range("J91").Value = 7752.5713
$ws.Range("L91").Value = 7752.5713
$ws.Range("N91").Value = -10560.5713
$ws.Range("H94").Value = 1468.7858
$ws.Range("I94").Value = 3765
$ws.Range("J94").Value = 550.3
$ws.Range("K94").Value = 3765
$ws.Range("L94").Value = 550.3
$ws.Range("M94").Value = -3314
$ws.Range("N94").Value = -1452.3
$ws.Range("H99").Value = 10035.204
$ws.Range("I99").Value = 6926.75
$ws.Range("J99").Value = 10725.973
$ws.Range("K99").Value = 6926.75
$ws.Range("L99").Value = 10725.973
$ws.Range("M99").Value = -5428.75
$ws.Range("N99").Value = -13721.973
$ws.Range("H105").Value = 1082.36
$ws.Range("I105").Value = 926.7619
$ws.Range("K105").Value = 926.7619
$ws.Range("M105").Value = 820.2381
$ws.Range("H107").Value = 1121.5714
$ws.Range("I107").Value = 1127.15
$ws.Range("J107").Value = 1107.625
$ws.Range("K107").Value = 1127.15
$ws.Range("L107").Value = 1107.625
$ws.Range("M107").Value = 792.8499999999999
$ws.Range("N107").Value = -4947.625
$ws.Range("H126").Value = 10035.204
$ws.Range("I126").Value = 6926.75
$ws.Range("J126").Value = 10725.973
$ws.Range("K126").Value = 20780.25
$ws.Range("L126").Value = 32177.919
$ws.Range("M126").Value = -18310.25
$ws.Range("N126").Value = -37117.919
$ws.Range("H132").Value = 22908.742
$ws.Range("I132").Value = 13246.05
$ws.Range("K132").Value = 39738.14999999999
$ws.Range("M132").Value = -37208.14999999999
$ws.Range("H141").Value = 198405.3
$ws.Range("J141").Value = 204555.62
$ws.Range("L141").Value = 204555.62
$ws.Range("N141").Value = -214915.62
$ws.Range("M4").Value = -49888

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 919.9
$ws.Range("J20").Value = 919.9
$ws.Range("L20").Value = 2759.7
$ws.Range("N20").Value = -3213.7
$ws.Range("H34").Value = 3757.8262
$ws.Range("J34").Value = 6928.0835
$ws.Range("L34").Value = 20784.2505
$ws.Range("N34").Value = -20952.2505
$ws.Range("H60").Value = 2923.2
$ws.Range("I60").Value = 141.33333
$ws.Range("J60").Value = 4777.778
$ws.Range("K60").Value = 423.99999
$ws.Range("L60").Value = 14333.334
$ws.Range("M60").Value = -172.99999
$ws.Range("N60").Value = -14835.334
$ws.Range("H75").Value = 2400
$ws.Range("J75").Value = 3600
$ws.Range("L75").Value = 10800
$ws.Range("N75").Value = -12796
$ws.Range("H78").Value = 2400
$ws.Range("J78").Value = 3600
$ws.Range("L78").Value = 32400
$ws.Range("N78").Value = -42384
$ws.Range("H103").Value = 753.44446
$ws.Range("I103").Value = 871
$ws.Range("J103").Value = 518.3333
$ws.Range("K103").Value = 2613
$ws.Range("L103").Value = 1554.9999
$ws.Range("M103").Value = -1734
$ws.Range("N103").Value = -3312.9999
$ws.Range("H104").Value = 10926.947
$ws.Range("I104").Value = 6874.143
$ws.Range("J104").Value = 12246.465
$ws.Range("K104").Value = 20622.429
$ws.Range("L104").Value = 36739.395
$ws.Range("M104").Value = -18001.429
$ws.Range("N104").Value = -41981.395
$ws.Range("H113").Value = 1447.95
$ws.Range("J113").Value = 1240
$ws.Range("L113").Value = 3720
$ws.Range("N113").Value = -8060
$ws.Range("H114").Value = 358.84616
$ws.Range("I114").Value = 416.83334
$ws.Range("J114").Value = 309.14285
$ws.Range("K114").Value = 1250.50002
$ws.Range("L114").Value = 927.4285500000001
$ws.Range("M114").Value = 2003.49998
$ws.Range("N114").Value = -7435.428550000001
$ws.Range("H120").Value = 9219.444
$ws.Range("I120").Value = 8797.200000000001
$ws.Range("J120").Value = 9747.25
$ws.Range("K120").Value = 26391.6
$ws.Range("L120").Value = 29241.75
$ws.Range("M120").Value = -21553.6
$ws.Range("N120").Value = -38917.75
$ws.Range("H129").Value = 2122.762
$ws.Range("I129").Value = 716.1111
$ws.Range("J129").Value = 3177.75
$ws.Range("K129").Value = 2148.3333
$ws.Range("L129").Value = 9533.25
$ws.Range("M129").Value = 2851.6667
$ws.Range("N129").Value = -19533.25
$ws.Range("H132").Value = 2086067.5
$ws.Range("I132").Value = 2255.1
$ws.Range("J132").Value = 5559088.5
$ws.Range("K132").Value = 20295.9
$ws.Range("L132").Value = 50031796.5
$ws.Range("M132").Value = -17765.9
$ws.Range("N132").Value = -50036856.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 500085
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 500085
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H69").Value = 160067
$ws.Range("J69").Value = 160067
$ws.Range("L69").Value = 160067
$ws.Range("H70").Value = 8816.115
$ws.Range("I70").Value = 7924.4707
$ws.Range("J70").Value = 10500.333
$ws.Range("K70").Value = 7924.4707
$ws.Range("L70").Value = 10500.333
$ws.Range("M70").Value = -7654.4707
$ws.Range("N70").Value = -11040.333
$ws.Range("H72").Value = 160067
$ws.Range("J72").Value = 160067
$ws.Range("L72").Value = 480201
$ws.Range("H73").Value = 8816.115
$ws.Range("I73").Value = 7924.4707
$ws.Range("J73").Value = 10500.333
$ws.Range("K73").Value = 7924.4707
$ws.Range("L73").Value = 10500.333
$ws.Range("M73").Value = -6988.4707
$ws.Range("N73").Value = -12372.333
$ws.Range("H80").Value = 1940
$ws.Range("I80").Value = 1946.25
$ws.Range("J80").Value = 1931.6666
$ws.Range("K80").Value = 1946.25
$ws.Range("L80").Value = 1931.6666
$ws.Range("M80").Value = -948.25
$ws.Range("N80").Value = -3927.6666
$ws.Range("H83").Value = 1940
$ws.Range("I83").Value = 1946.25
$ws.Range("J83").Value = 1931.6666
$ws.Range("K83").Value = 9731.25
$ws.Range("L83").Value = 9658.333000000001
$ws.Range("M83").Value = -4739.25
$ws.Range("N83").Value = -19642.333
$ws.Range("H98").Value = 29069
$ws.Range("J98").Value = 29069
$ws.Range("L98").Value = 29069
$ws.Range("H102").Value = 4870.6665
$ws.Range("I102").Value = 6799
$ws.Range("J102").Value = 1014
$ws.Range("K102").Value = 6799
$ws.Range("L102").Value = 1014
$ws.Range("M102").Value = -5177
$ws.Range("N102").Value = -4258
$ws.Range("H104").Value = 52500
$ws.Range("J104").Value = 85000
$ws.Range("L104").Value = 85000
$ws.Range("N104").Value = -91988
$ws.Range("H107").Value = 359.25
$ws.Range("J107").Value = 500
$ws.Range("L107").Value = 500
$ws.Range("H113").Value = 128018.125
$ws.Range("I113").Value = 157195
$ws.Range("J113").Value = 1585
$ws.Range("K113").Value = 157195
$ws.Range("L113").Value = 1585
$ws.Range("M113").Value = -155025
$ws.Range("N113").Value = -5925
$ws.Range("H126").Value = 9599.799999999999
$ws.Range("J126").Value = 9999.75
$ws.Range("L126").Value = 29999.25
$ws.Range("N126").Value = -34939.25
$ws.Range("H132").Value = 4441.636
$ws.Range("I132").Value = 3115.8235
$ws.Range("K132").Value = 9347.470499999999
$ws.Range("M132").Value = -6817.470499999999
$ws.Range("N69").Value = -161565
$ws.Range("N72").Value = -487689
$ws.Range("N98").Value = -35059
$ws.Range("N107").Value = -4340
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6894.143
$ws.Range("I7").Value = 6320.1816
$ws.Range("K7").Value = 6320.1816
$ws.Range("M7").Value = -6208.1816
$ws.Range("H16").Value = 8738
$ws.Range("I16").Value = 808.9375
$ws.Range("J16").Value = 24596.125
$ws.Range("K16").Value = 808.9375
$ws.Range("L16").Value = 24596.125
$ws.Range("M16").Value = -638.9375
$ws.Range("N16").Value = -24936.125
$ws.Range("H40").Value = 6180.25
$ws.Range("I40").Value = 6180.25
$ws.Range("K40").Value = 6180.25
$ws.Range("M40").Value = -6044.25
$ws.Range("H42").Value = 44974.332
$ws.Range("I42").Value = 44974.332
$ws.Range("K42").Value = 44974.332
$ws.Range("M42").Value = -44411.332
$ws.Range("H46").Value = 1510.4736
$ws.Range("I46").Value = 982.7
$ws.Range("K46").Value = 982.7
$ws.Range("M46").Value = -794.7
$ws.Range("H49").Value = 44974.332
$ws.Range("I49").Value = 44974.332
$ws.Range("K49").Value = 44974.332
$ws.Range("M49").Value = -44827.332
$ws.Range("H61").Value = 1260.6428
$ws.Range("I61").Value = 1061.1818
$ws.Range("J61").Value = 1992
$ws.Range("K61").Value = 1061.1818
$ws.Range("L61").Value = 1992
$ws.Range("M61").Value = -859.1818000000001
$ws.Range("N61").Value = -2396
$ws.Range("H68").Value = 3945.7693
$ws.Range("I68").Value = 2642.2856
$ws.Range("J68").Value = 5466.5
$ws.Range("K68").Value = 2642.2856
$ws.Range("L68").Value = 5466.5
$ws.Range("M68").Value = -1893.2856
$ws.Range("N68").Value = -6964.5
$ws.Range("H71").Value = 3945.7693
$ws.Range("I71").Value = 2642.2856
$ws.Range("J71").Value = 5466.5
$ws.Range("K71").Value = 13211.428
$ws.Range("L71").Value = 27332.5
$ws.Range("M71").Value = -9467.428
$ws.Range("N71").Value = -34820.5
$ws.Range("H94").Value = 69999.5
$ws.Range("J94").Value = 69999.5
$ws.Range("L94").Value = 69999.5
$ws.Range("N94").Value = -71351.5
$ws.Range("H98").Value = 30118.334
$ws.Range("J98").Value = 30118.334
$ws.Range("L98").Value = 30118.334
$ws.Range("N98").Value = -36108.334
$ws.Range("H99").Value = 34910.273
$ws.Range("I99").Value = 34910.273
$ws.Range("K99").Value = 34910.273
$ws.Range("M99").Value = -31915.273
$ws.Range("H113").Value = 1260.6428
$ws.Range("I113").Value = 1061.1818
$ws.Range("J113").Value = 1992
$ws.Range("K113").Value = 1061.1818
$ws.Range("L113").Value = 1992
$ws.Range("M113").Value = 1108.8182
$ws.Range("N113").Value = -6332
$ws.Range("H126").Value = 6894.143
$ws.Range("I126").Value = 6320.1816
$ws.Range("K126").Value = 18960.5448
$ws.Range("M126").Value = -16490.5448
$ws.Range("H132").Value = 7020.391
$ws.Range("I132").Value = 7124
$ws.Range("K132").Value = 21372
$ws.Range("M132").Value = -18842
$ws.Range("H136").Value = 4835.807
$ws.Range("I136").Value = 4259.479
$ws.Range("J136").Value = 7909.5557
$ws.Range("K136").Value = 12778.437
$ws.Range("L136").Value = 23728.6671
$ws.Range("M136").Value = -10228.437
$ws.Range("N136").Value = -28828.6671

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 58000
$ws.Range("J69").Value = 58000
$ws.Range("L69").Value = 58000
$ws.Range("N69").Value = -59498
$ws.Range("H72").Value = 58000
$ws.Range("J72").Value = 58000
$ws.Range("L72").Value = 174000
$ws.Range("N72").Value = -181488
$ws.Range("H82").Value = 49994.5
$ws.Range("J82").Value = 49994.5
$ws.Range("L82").Value = 49994.5
$ws.Range("N82").Value = -50760.5
$ws.Range("H85").Value = 49994.5
$ws.Range("J85").Value = 49994.5
$ws.Range("L85").Value = 49994.5
$ws.Range("N85").Value = -52646.5
$ws.Range("H105").Value = 34750
$ws.Range("J105").Value = 34750
$ws.Range("L105").Value = 34750
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H113").Value = 702.1429000000001
$ws.Range("I113").Value = 559.25
$ws.Range("K113").Value = 1677.75
$ws.Range("M113").Value = 492.25
$ws.Range("H122").Value = 6209.143
$ws.Range("I122").Value = 4293
$ws.Range("K122").Value = 12879
$ws.Range("M122").Value = -10429
$ws.Range("H136").Value = 10716702
$ws.Range("I136").Value = 13045837
$ws.Range("K136").Value = 39137511
$ws.Range("M136").Value = -39134961
$ws.Range("N105").Value = -41738
$ws.Range("N110").ClearContents()

Write-Host "Applied all updates."